$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "64.195.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.135.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -1.23%  "
$ws.Range("E4").Value2 = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "570.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "161.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -4.39%  "
$ws.Range("E7").Value2 = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -6.97%  "
$ws.Range("E9").Value2 = "  -3.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "6.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.381"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "3.682.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -1.16%  "
$ws.Range("E13").Value2 = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "64.268.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "24.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "3.133.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -1.12%  "
$ws.Range("E17").Value2 = "  -3.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "399.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -4.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "5.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -2.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "67.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -3.44%  "
$ws.Range("E25").Value2 = "  -1.42%  "
$ws.Range("E26").Value2 = "  -5.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.0₂01000"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -5.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "8.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.988"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -1.30%  "
$ws.Range("E30").Value2 = "  +0.11%  "
$ws.Range("E31").Value2 = "  -1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "21.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -2.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "158.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +0.85%  "
$ws.Range("E34").Value2 = "  -1.56%  "
$ws.Range("E35").Value2 = "  -4.99%  "
$ws.Range("E36").Value2 = "  -2.89%  "
$ws.Range("E37").Value2 = "  -2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.656.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -2.44%  "
$ws.Range("E39").Value2 = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "23.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -3.04%  "
$ws.Range("E41").Value2 = "  -2.59%  "
$ws.Range("E42").Value2 = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.687"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -3.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.0609"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "5.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -3.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "288.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -1.57%  "
$ws.Range("E47").Value2 = "  -3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "20.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -2.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0971"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -1.79%  "
$ws.Range("E51").Value2 = "  -0.02%  "
